$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 13 - Geri Vass: add Walk-On Music "Jayron, GEWOONRAVES - Knock Knock"
$ws.Range("C13").Value = "Jayron, GEWOONRAVES - Knock Knock"

# Row 14 - Balazs Papai: add Nickname "Baja"
$ws.Range("B14").Value = "Baja"

# Row 13 - Geri Vass: add Nickname "Papito" (match formatting of the rest of the row)
$ws.Range("E13").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B13").Value = "Papito"

# Update the active selection to B16 (was C16)
$ws.Range("B16").Select()
